# Update the "dSF" (column F) values for the rows that were recalculated
# after repulling/pushing data (per commit message: "repull data, push
# all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    6  = 3
    9  = 1
    11 = -1
    12 = -1
    13 = -2
    16 = -2
    19 = 2
    22 = 0
    24 = 0
    25 = -1
    26 = -1
    29 = 1
    30 = 0
    33 = 2
    36 = 2
    40 = 3
    41 = 0
    44 = -2
    56 = 2
    64 = 0
    65 = -1
    71 = -2
    73 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
